$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. URL value changed (row 2, column B): pythia -> cicada
$ws.Cells.Item(2,2).Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/EvalStatus"

# 2. Date value changed (row 8, column B)
$ws.Cells.Item(8,2).Value = "2026-02-11T14:37:07-05:00"

# 3. Insert a new "Jurisdiction" row (with an empty value) right after the
#    "Contact" row (row 10) and before "Description" (row 11), pushing every
#    row from 11..21 down by one (to 12..22). Shift bottom-up, using
#    Copy(destination) so the string-vs-number cell type and the cell style
#    (s="2") are preserved exactly - a plain Value round-trip or Rows.Insert
#    would re-type "4" as a number or mint a fresh, unused style.
for ($r = 21; $r -ge 11; $r--) {
    $ws.Cells.Item($r+1,1).ClearContents()
    $ws.Cells.Item($r+1,2).ClearContents()
    $ws.Cells.Item($r,1).Copy($ws.Cells.Item($r+1,1))
    $ws.Cells.Item($r,2).Copy($ws.Cells.Item($r+1,2))
}

# Row 11 is now free - make it the new Jurisdiction row (Value column blank)
$ws.Cells.Item(11,2).ClearContents()
$ws.Cells.Item(11,1).Value = "Jurisdiction"
